# Applies the edits described by the target diff:
#  1. Update the auto date placeholder text ("datetimeFigureOut" field) from
#     12/23/2025 to 1/9/2026 on the slide master and on every slide layout.
#  2. On slide 3, split the title run so the phrase "YÊU CẦU " is removed
#     from "MÔ TẢ ĐỀ TÀI & YÊU CẦU NGHIỆP VỤ " leaving
#     "MÔ TẢ ĐỀ TÀI & NGHIỆP VỤ ".
#  3. On slide 3, remove the leading "YÊU CẦU " from the requirements
#     textbox heading "YÊU CẦU NGHIỆP VỤ " leaving "NGHIỆP VỤ ".

$p = $ppt.ActivePresentation

# --- 1. Date placeholder (footer "date" field) -----------------------------
function Set-DatePlaceholderText {
    param($container, $newText)
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

Set-DatePlaceholderText $p.SlideMaster "1/9/2026"
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    Set-DatePlaceholderText $p.SlideMaster.CustomLayouts.Item($L) "1/9/2026"
}

# --- 2 & 3. Slide 3 text edits ---------------------------------------------
$slide3 = $p.Slides.Item(3)

# Title shape: "MÔ TẢ ĐỀ TÀI & YÊU CẦU NGHIỆP VỤ " -> "MÔ TẢ ĐỀ TÀI & NGHIỆP VỤ "
$titleShape = $slide3.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "MÔ TẢ ĐỀ TÀI & NGHIỆP VỤ "
$titleMiddle = $titleRange.Characters(10, 6)
$titleMiddle.Text = "TÀI & "

# Requirements textbox: "YÊU CẦU NGHIỆP VỤ " -> "NGHIỆP VỤ "
$reqShape = $slide3.Shapes.Item(4)
$reqHeading = $reqShape.TextFrame.TextRange.Characters(1, 8)
$reqHeading.Text = ""
